$wb = $excel.ActiveWorkbook

# Rename the active sheet from "Sheet1_2(19Nov)" to "Sheet1_2(20Nov)"
$ws = $wb.ActiveSheet
$ws.Name = "Sheet1_2(20Nov)"

# Fix the typo in K13 (38863 -> 3863)
$ws.Range("K13").Value = 3863

# Fill in the new "day 6" (column L) readings
$ws.Range("L8").Value = 370
$ws.Range("L9").Value = 537
$ws.Range("L10").Value = 182
$ws.Range("L11").Value = 578
$ws.Range("L12").Value = 829
$ws.Range("L13").Value = 3861
$ws.Range("L14").Value = 992
$ws.Range("L15").Value = 10

# Update the selected cell in the sheet view
$ws.Range("L15").Select()

$excel.Calculate()
